# Add a new "DevHub" app entry for user j_thomas on the Attributes sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes")

# Insert a new row right after the existing "j_thomas | app | QMC" row (row 56),
# shifting the "udc" block (rows 57-66) down by one.
[void]$ws.Rows.Item(57).Insert()

$ws.Range("A57").Value = "j_thomas"
$ws.Range("B57").Value = "app"
$ws.Range("C57").Value = "DevHub"

[void]$ws.Range("C57").Select()
